# Auto update Excel log
# Appends newly-logged sensor rows to the "PIR" sheet (rows 94-106) and the
# "Humidity" sheet (rows 65-69), matching the trailing tail of data already
# present on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PIR sheet: append rows 94-106 (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-01-30", "18:26:23", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:26:24", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:26:28", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:26:33", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:26:38", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:26:43", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:26:48", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:26:53", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:26:58", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:27:03", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:27:08", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:27:13", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:27:18", "18:00", "Bathroom", "No Motion", "Inactive")
)

$startRow = 94
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $values = $pirRows[$i]

    # Column A holds a date-shaped string ("2026-01-30"). Excel would
    # normally auto-convert that to a real date serial, so force the cell
    # to Text first, write the literal string, then drop the number
    # format back to Normal so no stray style lingers on the cell.
    $cellA = $pir.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.Style = "Normal"

    $pir.Cells.Item($r, 2).Value = $values[1]
    $pir.Cells.Item($r, 3).Value = $values[2]
    $pir.Cells.Item($r, 4).Value = $values[3]
    $pir.Cells.Item($r, 5).Value = $values[4]
    $pir.Cells.Item($r, 6).Value = $values[5]
}

# ---------------------------------------------------------------------
# Humidity sheet: append rows 65-69 (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")

$humidityRows = @(
    @("2026-01-30", "18:26:24", "18:00", "Bathroom", "86.6%", "Active"),
    @("2026-01-30", "18:26:28", "18:00", "Bathroom", "86.6%", "Active"),
    @("2026-01-30", "18:26:49", "18:00", "Bathroom", "86.7%", "Active"),
    @("2026-01-30", "18:27:09", "18:00", "Bathroom", "86.6%", "Active"),
    @("2026-01-30", "18:27:14", "18:00", "Bathroom", "86.6%", "Active")
)

$startRow = 65
for ($i = 0; $i -lt $humidityRows.Count; $i++) {
    $r = $startRow + $i
    $values = $humidityRows[$i]

    $cellA = $humidity.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.Style = "Normal"

    $humidity.Cells.Item($r, 2).Value = $values[1]
    $humidity.Cells.Item($r, 3).Value = $values[2]
    $humidity.Cells.Item($r, 4).Value = $values[3]

    # Column E holds a percentage-shaped string ("86.6%"). Same
    # auto-conversion problem as the date column above, so apply the
    # same Text-format-then-reset trick.
    $cellE = $humidity.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $values[4]
    $cellE.Style = "Normal"

    $humidity.Cells.Item($r, 6).Value = $values[5]
}
